$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force column D to Text format first so numeric-looking price strings
# (e.g. '596.49', '1.00') are stored as text, matching the source inlineStr cells,
# instead of being auto-coerced to numbers by Excel's smart entry.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.381.81"
$ws.Range("E2").Value = "  -1.55%  "

$ws.Range("D3").Value = "2.671.63"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "596.49"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").Value = "162.35"
$ws.Range("E6").Value = "  +2.66%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").Value = "2.672.00"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("E11").Value = "  +0.50%  "

$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").Value = "5.19"
$ws.Range("E13").Value = "  -1.50%  "

$ws.Range("D14").Value = "27.72"
$ws.Range("E14").Value = "  -1.24%  "

$ws.Range("D15").Value = "3.166.40"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").Value = "0.0000183"
$ws.Range("E16").Value = "  -1.82%  "

$ws.Range("D17").Value = "67.391.33"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "2.675.44"
$ws.Range("E18").Value = "  +1.58%  "

$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  -1.47%  "

$ws.Range("D20").Value = "362.82"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").Value = "7.46"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").Value = "4.37"
$ws.Range("E22").Value = "  -2.66%  "

$ws.Range("D23").Value = "4.80"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("E24").Value = "  -3.69%  "

$ws.Range("E25").Value = "  -3.98%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").Value = "2.828.90"

$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").Value = "552.74"
$ws.Range("E31").Value = "  -4.21%  "

$ws.Range("D32").Value = "7.92"
$ws.Range("E32").Value = "  -3.79%  "

$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("D34").Value = "1.93"
$ws.Range("E34").Value = "  +1.39%  "

$ws.Range("E35").Value = "  -0.77%  "

$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").Value = "1.57"
$ws.Range("E37").Value = "  -4.85%  "

$ws.Range("D38").Value = "19.49"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").Value = "156.68"
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  -1.75%  "

$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  -3.34%  "

$ws.Range("D42").Value = "5.25"
$ws.Range("E42").Value = "  -1.72%  "

$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("D46").Value = "40.34"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "0.0₆0301"
$ws.Range("E47").Value = "  -4.84%  "

$ws.Range("D48").Value = "0.585"
$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("D49").Value = "152.88"
$ws.Range("E49").Value = "  -3.72%  "

$ws.Range("D50").Value = "3.82"
$ws.Range("E50").Value = "  -2.93%  "

$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  -0.76%  "

# Strip the temporary Text number-format back off so the cell styling
# matches the original workbook (which used the default/general style).
$ws.Range("D2:D51").ClearFormats()
